$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '60.192.45'
$ws.Range("E2").Value = '  +3.01%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.569.80'
$ws.Range("E3").Value = '  +4.56%  '

$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '505.04'
$ws.Range("E5").Value = '  +1.87%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '153.27'
$ws.Range("E6").Value = '  -3.79%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  +0.30%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.578'
$ws.Range("E8").Value = '  -6.15%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.588.17'
$ws.Range("E9").Value = '  +4.24%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.64'
$ws.Range("E10").Value = '  +6.05%  '

$ws.Range("E11").Value = '  +0.51%  '

$ws.Range("E12").Value = '  +1.33%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.027.78'
$ws.Range("E14").Value = '  +5.08%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '60.242.61'
$ws.Range("E15").Value = '  +3.26%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.55'
$ws.Range("E16").Value = '  +0.60%  '

$ws.Range("E17").Value = '  +2.85%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.585.31'
$ws.Range("E18").Value = '  +3.83%  '

$ws.Range("E19").Value = '  +1.18%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '345.82'
$ws.Range("E20").Value = '  +5.06%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.27'
$ws.Range("E21").Value = '  +1.12%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.04'
$ws.Range("E22").Value = '  +1.07%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.12%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '60.22'
$ws.Range("E24").Value = '  +2.61%  '

$ws.Range("E25").Value = '  +2.14%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.165'
$ws.Range("E26").Value = '  +0.88%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.688.52'
$ws.Range("E27").Value = '  +4.67%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.995'

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '0.0₃0848'
$ws.Range("E29").Value = '  +4.96%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '7.44'
$ws.Range("E30").Value = '  +0.13%  '

$ws.Range("E31").Value = '  +0.12%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '155.14'
$ws.Range("E32").Value = '  +2.45%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.16'
$ws.Range("E33").Value = '  +0.68%  '

$ws.Range("E34").Value = '  +1.03%  '

$ws.Range("E36").Value = '  +3.36%  '

$ws.Range("E37").Value = '  +1.46%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.857'
$ws.Range("E38").Value = '  +21.56%  '

$ws.Range("B39").Value = 'Filecoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.77'
$ws.Range("E39").Value = '  +3.35%  '

$ws.Range("B40").Value = 'Fetch.AI'
$ws.Range("C40").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.841'
$ws.Range("E40").Value = '  -0.86%  '

$ws.Range("E41").Value = '  +3.22%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '299.92'
$ws.Range("E42").Value = '  +5.40%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '35.48'
$ws.Range("E43").Value = '  +3.13%  '

$ws.Range("E44").Value = '  +3.20%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.615'
$ws.Range("E45").Value = '  +1.10%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0993'
$ws.Range("E46").Value = '  -1.69%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.996'
$ws.Range("E47").Value = '  +0.47%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '19.69'
$ws.Range("E48").Value = '  +8.58%  '

$ws.Range("E49").Value = '  +3.24%  '

$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0233'
$ws.Range("E50").Value = '  -0.99%  '

$ws.Range("B51").Value = 'Maker'
$ws.Range("C51").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.014.25'
$ws.Range("E51").Value = '  +5.38%  '
